$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.397.19"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "1.846.36"
$ws.Range("E3").Value = "  -0.28%  "

$style_D4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("D4").Style = $style_D4
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("E5").Value = "  -0.03%  "

$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6299"
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("E7").Value = "  -0.17%  "

$style_D8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07486"
$ws.Range("D8").Style = $style_D8
$ws.Range("E8").Value = "  -2.19%  "

$style_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2907"
$ws.Range("D9").Style = $style_D9
$ws.Range("E9").Value = "  -0.17%  "

$style_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.40"
$ws.Range("D10").Style = $style_D10
$ws.Range("E10").Value = "  -1.75%  "

$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").Value = "1.846.57"
$ws.Range("E12").Value = "  -2.21%  "

$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.004"
$ws.Range("D13").Style = $style_D13
$ws.Range("E13").Value = "  -0.78%  "

$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6784"
$ws.Range("D14").Style = $style_D14
$ws.Range("E14").Value = "  -0.56%  "

$style_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.13"
$ws.Range("D16").Style = $style_D16
$ws.Range("E16").Value = "  -1.60%  "

$style_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.140"
$ws.Range("D17").Style = $style_D17
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("D18").Value = "29.415.74"
$ws.Range("E18").Value = "  -0.38%  "

$style_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.21"
$ws.Range("D19").Style = $style_D19
$ws.Range("E19").Value = "  -0.22%  "

$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.31"
$ws.Range("D20").Style = $style_D20
$ws.Range("E20").Value = "  -0.31%  "

$style_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("D21").Style = $style_D21
$ws.Range("E21").Value = "  -0.21%  "

$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.436"
$ws.Range("D22").Style = $style_D22
$ws.Range("E22").Value = "  -0.44%  "

$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9996"
$ws.Range("D23").Style = $style_D23
$ws.Range("E23").Value = "  -0.25%  "

$style_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.98"
$ws.Range("D24").Style = $style_D24
$ws.Range("E24").Value = "  +0.54%  "

$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1375"
$ws.Range("D25").Style = $style_D25
$ws.Range("E25").Value = "  -0.63%  "

$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.421"
$ws.Range("D26").Style = $style_D26
$ws.Range("E26").Value = "  -0.16%  "

$style_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.55"
$ws.Range("D27").Style = $style_D27
$ws.Range("E27").Value = "  -0.99%  "

$style_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06328"
$ws.Range("D28").Style = $style_D28
$ws.Range("E28").Value = "  +12.87%  "

$style_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.392"
$ws.Range("D29").Style = $style_D29
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  +0.48%  "

$style_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.093"
$ws.Range("D31").Style = $style_D31
$ws.Range("E31").Value = "  -1.10%  "

$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.056"
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("E33").Value = "  -1.25%  "

$ws.Range("E34").Value = "  -2.01%  "

$style_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6931"
$ws.Range("D35").Style = $style_D35
$ws.Range("E35").Value = "  -0.50%  "

$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.581"
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("D37").Value = "1.258.18"
$ws.Range("E37").Value = "  +2.35%  "

$ws.Range("E38").Value = "  +3.94%  "

$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01819"
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = "  +0.71%  "

$style_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.549"
$ws.Range("D40").Style = $style_D40
$ws.Range("E40").Value = "  +1.63%  "

$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9090"
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = "  +0.22%  "

$style_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9995"
$ws.Range("D42").Style = $style_D42
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").Value = "2.007.70"
$ws.Range("E43").Value = "  -14.63%  "

$style_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.33"
$ws.Range("D44").Style = $style_D44
$ws.Range("E44").Value = "  -1.09%  "

$style_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.33"
$ws.Range("D45").Style = $style_D45
$ws.Range("E45").Value = "  +0.39%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$style_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.055"
$ws.Range("D46").Style = $style_D46
$ws.Range("E46").Value = "  -2.08%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1173"
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = "  +1.73%  "

$ws.Range("E48").Value = "  -2.25%  "

$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.040"
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = "  +0.01%  "

$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.682"
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = "  -0.07%  "

$style_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3939"
$ws.Range("D51").Style = $style_D51
$ws.Range("E51").Value = "  -2.18%  "
